$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "IMF (20%)" values that will occupy columns F/G
$imf20Sales = @{
    2 = 0.5001387244282065
    3 = 0.5499393467648322
    4 = 0.1011494731738439
    5 = 0.2211979425309333
    6 = 0.8501823953898133
    7 = 0.4341752803803408
    8 = 0.4123631444047267
    9 = 0.7467061426163347
    10 = 0.4552319623762475
    11 = 1.486349685879603
}
$imf20SalesEmp = @{
    2 = 0.4034859978630822
    3 = 0.629475953570518
    4 = 0.4811403687082352
    5 = 0.7602146625722329
    6 = 0.6457630892642331
    7 = 1.501722643141881
    8 = 0.4537896841353516
    9 = 0.6101153731432726
    10 = 0.4861428706789624
    11 = 0.9247074269381814
}

for ($row = 2; $row -le 11; $row++) {
    # Capture the current F/G values (old "IMF - Sales" / "IMF - Sales + Emp") - they move to H/I
    $oldF = $ws.Cells.Item($row, 6).Value2
    $oldG = $ws.Cells.Item($row, 7).Value2

    # H/I now hold what used to be in F/G
    $ws.Cells.Item($row, 8).Value = $oldF
    $ws.Cells.Item($row, 9).Value = $oldG

    # F/G now hold the new "IMF (20%)" values
    $ws.Cells.Item($row, 6).Value = $imf20Sales[$row]
    $ws.Cells.Item($row, 7).Value = $imf20SalesEmp[$row]
}

# Update header labels (row 1) to reflect the new shared-string text/order
$ws.Cells.Item(1, 6).Value = "IMF (20%) - Sales"
$ws.Cells.Item(1, 7).Value = "IMF (20%) - Sales + Emp"
$ws.Cells.Item(1, 8).Value = "IMF - Sales"
$ws.Cells.Item(1, 9).Value = "IMF - Sales + Emp"
